$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.416940450668335
$ws.Range("B1").Value = 2.250223398208618
$ws.Range("C1").Value = 2.51053786277771
$ws.Range("D1").Value = 3.083885669708252
$ws.Range("E1").Value = 0.9102909564971924
